# Add season-record columns (Wins / Losses / Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — same header style as the other header cells (A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-50 — constant season record for every player row.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 90  # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 72  # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF -> Ties
}

Write-Output "season record columns added"
